# The post formerly stored at row 263 ("「お気に入りのおもちゃには触らせない」...")
# was removed from the source data. Delete that entire row; Excel shifts every
# row below it up by one (rows 264-314 become 263-313), and the sheet's used
# range shrinks from A1:C314 to A1:C313 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(263).Delete()
